# Auto update Excel log
# Appends newly-logged sensor readings (2026-01-28, ~12:33-12:34) to the
# PIR, Humidity, Temperature and Proximity sheets of the SeniorConnect
# master log.
#
# Values such as "2026-01-28" or "87.3%" would normally be auto-coerced by
# Excel into a date/percentage number on assignment. Every row in this log
# is stored as literal text instead (matching all pre-existing rows), so
# each cell is written with a leading single-quote (forces text entry,
# exactly like typing '2026-01-28 into Excel) and then has its style reset
# to match the plain, unformatted style already used throughout the sheet
# (clears the transient "quote prefixed" style Excel applies).

$wb = $excel.ActiveWorkbook

function Set-LogCell($sheet, $row, $col, $val) {
    $sheet.Cells.Item($row, $col).Value = "'" + $val
    $sheet.Cells.Item($row, $col).Style = $sheet.Cells.Item(1, 1).Style
}

function Add-LogRow($sheet, $rowData) {
    # $rowData = @(row, Date, Timestamp, Hour, Location, Value, Status)
    Set-LogCell $sheet $rowData[0] 1 $rowData[1]
    Set-LogCell $sheet $rowData[0] 2 $rowData[2]
    Set-LogCell $sheet $rowData[0] 3 $rowData[3]
    Set-LogCell $sheet $rowData[0] 4 $rowData[4]
    Set-LogCell $sheet $rowData[0] 5 $rowData[5]
    Set-LogCell $sheet $rowData[0] 6 $rowData[6]
}

# ---------------------------------------------------------------------
# PIR sheet: append rows 347-359
# ---------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @(347, "2026-01-28", "12:33:21", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(348, "2026-01-28", "12:33:24", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(349, "2026-01-28", "12:33:28", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(350, "2026-01-28", "12:33:35", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(351, "2026-01-28", "12:33:39", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(352, "2026-01-28", "12:33:43", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(353, "2026-01-28", "12:33:49", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(354, "2026-01-28", "12:33:54", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(355, "2026-01-28", "12:33:59", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(356, "2026-01-28", "12:34:04", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(357, "2026-01-28", "12:34:09", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(358, "2026-01-28", "12:34:15", "12:00", "Bathroom", "No Motion", "Inactive"),
    @(359, "2026-01-28", "12:34:19", "12:00", "Bathroom", "No Motion", "Inactive")
)
foreach ($r in $pirRows) {
    Add-LogRow $wsPIR $r
}

# ---------------------------------------------------------------------
# Humidity sheet: append rows 323-334
# ---------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @(323, "2026-01-28", "12:33:22", "12:00", "Bathroom", "87.3%", "Active"),
    @(324, "2026-01-28", "12:33:25", "12:00", "Bathroom", "87.2%", "Active"),
    @(325, "2026-01-28", "12:33:33", "12:00", "Bathroom", "87.2%", "Active"),
    @(326, "2026-01-28", "12:33:37", "12:00", "Bathroom", "87.2%", "Active"),
    @(327, "2026-01-28", "12:33:41", "12:00", "Bathroom", "86.3%", "Active"),
    @(328, "2026-01-28", "12:33:45", "12:00", "Bathroom", "87.3%", "Active"),
    @(329, "2026-01-28", "12:33:57", "12:00", "Bathroom", "87.3%", "Active"),
    @(330, "2026-01-28", "12:34:01", "12:00", "Bathroom", "86.4%", "Active"),
    @(331, "2026-01-28", "12:34:05", "12:00", "Bathroom", "87.2%", "Active"),
    @(332, "2026-01-28", "12:34:10", "12:00", "Bathroom", "86.2%", "Active"),
    @(333, "2026-01-28", "12:34:13", "12:00", "Bathroom", "87.2%", "Active"),
    @(334, "2026-01-28", "12:34:17", "12:00", "Bathroom", "87.2%", "Active")
)
foreach ($r in $humidityRows) {
    Add-LogRow $wsHumidity $r
}

# ---------------------------------------------------------------------
# Temperature sheet: append rows 323-334
# ---------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @(323, "2026-01-28", "12:33:23", "12:00", "Bathroom", "23.1C", "Active"),
    @(324, "2026-01-28", "12:33:26", "12:00", "Bathroom", "23.0C", "Active"),
    @(325, "2026-01-28", "12:33:34", "12:00", "Bathroom", "23.0C", "Active"),
    @(326, "2026-01-28", "12:33:38", "12:00", "Bathroom", "23.0C", "Active"),
    @(327, "2026-01-28", "12:33:42", "12:00", "Bathroom", "23.0C", "Active"),
    @(328, "2026-01-28", "12:33:46", "12:00", "Bathroom", "23.0C", "Active"),
    @(329, "2026-01-28", "12:33:58", "12:00", "Bathroom", "23.0C", "Active"),
    @(330, "2026-01-28", "12:34:02", "12:00", "Bathroom", "23.0C", "Active"),
    @(331, "2026-01-28", "12:34:06", "12:00", "Bathroom", "23.0C", "Active"),
    @(332, "2026-01-28", "12:34:10", "12:00", "Bathroom", "23.0C", "Active"),
    @(333, "2026-01-28", "12:34:14", "12:00", "Bathroom", "23.0C", "Active"),
    @(334, "2026-01-28", "12:34:18", "12:00", "Bathroom", "23.0C", "Active")
)
foreach ($r in $temperatureRows) {
    Add-LogRow $wsTemperature $r
}

# ---------------------------------------------------------------------
# Proximity sheet: append row 7
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
Add-LogRow $wsProximity @(7, "2026-01-28", "12:34:21", "12:00", "Living Room Main Entrance", "Detected", "Active")

Write-Output "Appended $($pirRows.Count) PIR rows, $($humidityRows.Count) Humidity rows, $($temperatureRows.Count) Temperature rows, and 1 Proximity row."
